$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.940.98'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.61%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.889.89'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.03%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '249.86'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.25%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.84%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.02'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.349'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '51.37'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.21%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0740'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.21%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0974'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.163.80'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.82'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.95%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.717'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.19%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.91'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.892.77'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '35.003.60'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.22'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0825'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '247.70'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.98%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.82'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.90%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.95'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.28%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.26%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.81%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.24'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.32'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.42'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.25'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.95%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.128.35'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.61%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0583'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.28%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'TrustWalletToken'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.56'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +5.58%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.88'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.13%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'InternetComputer(DFINITY)'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.17'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.16%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'BinanceUSD'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.01'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.28%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -8.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.99'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.49%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.28'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.36%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.05'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.46%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0671'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.09%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.97%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.43%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.295.07'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -4.66%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.36'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.72%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.12%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.74'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.75%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '12.14'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0762'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +6.61%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.46'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.29%  '
